# Stand-up Meeting Sprint6.xlsx — update the team member's daily status
# entries in rows 17-19 with the new sprint's motives.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17 ("1)" items) ---
$ws.Range("B17").Value = "1)Learned about recycler view and components of the android"
$ws.Range("C17").Value = "1) Prepared slides for the workshop"
$ws.Range("D17").Value = "1)I prepared and given presentation today"

# --- Row 18 ("2)" items) ---
$ws.Range("B18").Value = "2)Preparing slides for presentation"
$ws.Range("C18").Value = "2)I will give presentation today"
$ws.Range("D18").Value = "2)Discussing about pros and cons heroku progress today"

# --- Row 19 ("3)" items) ---
$ws.Range("B19").Value = "3)NA"
$ws.Range("C19").Value = "3) NA"
$ws.Range("D19").Value = "3) Finalizing the database"

# The longer text in rows 17-18 now wraps onto two lines, so those rows
# grow from single-line (15.5) to double-line (31) height; row 19's text
# still fits on one line.
$ws.Rows("17:18").RowHeight = 31

# Update the view: scrolled down a bit further and a different cell selected.
$ws.Range("D18").Select()
